# Java Sockets understands byte array messages. Timing collected.
#
# Adds a new "JAVA SOCKETS RESULTS" worksheet (mirroring the existing
# "WINSOCK RESULTS" sheet layout) populated with the collected timings,
# and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# The existing "WINSOCK RESULTS" sheet (4th tab) is the template this
# new results sheet is modeled after - reuse its formatting.
$ws4 = $wb.Worksheets.Item(4)

# New sheet goes at the end of the tab strip, right after WINSOCK RESULTS.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws5.Name = "JAVA SOCKETS RESULTS"

# --- Formatting: copy over just the formats used on WINSOCK RESULTS ---
$ws4.Range("A1:E1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)

$ws4.Range("B2:E2").Copy()
$ws5.Range("B2").PasteSpecial(-4122)

$ws4.Range("A3:A12").Copy()
$ws5.Range("A3").PasteSpecial(-4122)

$ws4.Range("A13").Copy()
$ws5.Range("A13").PasteSpecial(-4122)

# --- Title row ---
$ws5.Range("A1").Value = "Java Sockets (No JNI)"
$ws5.Range("A1:E1").Merge()

# --- Column headers ---
$ws5.Range("B2").Value = "40 BYTES "
$ws5.Range("C2").Value = "400 BYTES"
$ws5.Range("D2").Value = "4000 BYTES"
$ws5.Range("E2").Value = "40 000 BYTES"

# --- Run numbers (column A) ---
$runs = New-Object 'object[,]' 10,1
for ($i = 0; $i -lt 10; $i++) { $runs[$i,0] = $i + 1 }
$ws5.Range("A3:A12").Value = $runs

# --- Timing data (columns B:E) ---
$data = New-Object 'object[,]' 10,4
$data[0,0]=5913293;  $data[0,1]=5980502;  $data[0,2]=7009409;  $data[0,3]=6315794
$data[1,0]=6116808;  $data[1,1]=6142107;  $data[1,2]=6151546;  $data[1,3]=6184396
$data[2,0]=6149657;  $data[2,1]=6368655;  $data[2,2]=6008065;  $data[2,3]=5937080
$data[3,0]=6090000;  $data[3,1]=6108124;  $data[3,2]=5711664;  $data[3,3]=6934271
$data[4,0]=6074142;  $data[4,1]=6122472;  $data[4,2]=6719805;  $data[4,3]=5744891
$data[5,0]=5859676;  $data[5,1]=5820786;  $data[5,2]=5751310;  $data[5,3]=6995817
$data[6,0]=6604642;  $data[6,1]=5925375;  $data[6,2]=6303711;  $data[6,3]=5621799
$data[7,0]=5592726;  $data[7,1]=5822673;  $data[7,2]=6430200;  $data[7,3]=5738095
$data[8,0]=5438673;  $data[8,1]=6219887;  $data[8,2]=5604431;  $data[8,3]=6887450
$data[9,0]=6552159;  $data[9,1]=5653895;  $data[9,2]=6486082;  $data[9,3]=6496654
$ws5.Range("B3:E12").Value = $data

# --- Averages row ---
$ws5.Range("A13").Value = "Average"
$ws5.Range("B13").Formula = "=AVERAGE(B3:B12)"
$ws5.Range("C13").Formula = "=AVERAGE(C3:C12)"
$ws5.Range("D13").Formula = "=AVERAGE(D3:D12)"
$ws5.Range("E13").Formula = "=AVERAGE(E3:E12)"

# --- Column widths (match the fitted widths used on the other sheets) ---
$ws5.Columns.Item(2).ColumnWidth = 8.25
$ws5.Columns.Item(3).ColumnWidth = 8.875
$ws5.Columns.Item(4).ColumnWidth = 9.875
$ws5.Columns.Item(5).ColumnWidth = 11.25

# --- Selections: WINSOCK RESULTS loses its tab-selected/cursor state,
#     new JAVA SOCKETS RESULTS sheet becomes the active tab ---
$ws4.Range("B2:E2").Select()
$ws5.Range("E13").Select()
$ws5.Activate()
